$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "23.192.12"
Set-TextValue $ws "E2" "  -3.42%  "
Set-TextValue $ws "D3" "1.607.54"
Set-TextValue $ws "E3" "  -2.94%  "
Set-TextValue $ws "E4" "  +0.14%  "
Set-TextValue $ws "D5" "1.000"
Set-TextValue $ws "E5" "  -0.01%  "
Set-TextValue $ws "D6" "302.76"
Set-TextValue $ws "E6" "  -2.30%  "
Set-TextValue $ws "D7" "0.3768"
Set-TextValue $ws "E7" "  -3.38%  "
Set-TextValue $ws "D8" "0.3649"
Set-TextValue $ws "E8" "  -5.04%  "
Set-TextValue $ws "D9" "48.99"
Set-TextValue $ws "E9" "  -4.69%  "
Set-TextValue $ws "D10" "1.000"
Set-TextValue $ws "E10" "  +0.10%  "
Set-TextValue $ws "D11" "1.270"
Set-TextValue $ws "E11" "  -6.42%  "
Set-TextValue $ws "D12" "0.08082"
Set-TextValue $ws "E12" "  -4.43%  "
Set-TextValue $ws "D13" "23.06"
Set-TextValue $ws "E13" "  -3.93%  "
Set-TextValue $ws "D14" "6.581"
Set-TextValue $ws "E14" "  -7.77%  "
Set-TextValue $ws "D15" "7.569"
Set-TextValue $ws "E15" "  -4.43%  "
Set-TextValue $ws "D16" "0.00001267"
Set-TextValue $ws "E16" "  -3.91%  "
Set-TextValue $ws "D17" "1.610.01"
Set-TextValue $ws "E17" "  -2.68%  "
Set-TextValue $ws "D18" "91.55"
Set-TextValue $ws "E18" "  -3.26%  "
Set-TextValue $ws "D19" "0.06776"
Set-TextValue $ws "E19" "  -3.18%  "
Set-TextValue $ws "D20" "18.35"
Set-TextValue $ws "E20" "  -7.36%  "
Set-TextValue $ws "D21" "6.579"
Set-TextValue $ws "E21" "  -5.26%  "
Set-TextValue $ws "D22" "1.0000"
Set-TextValue $ws "E22" "  -0.07%  "
Set-TextValue $ws "D23" "13.10"
Set-TextValue $ws "E23" "  -4.48%  "
Set-TextValue $ws "D24" "23.233.05"
Set-TextValue $ws "E24" "  -3.24%  "
Set-TextValue $ws "D25" "2.356"
Set-TextValue $ws "E25" "  -5.08%  "
Set-TextValue $ws "D26" "2.916"
Set-TextValue $ws "E26" "  -3.09%  "
Set-TextValue $ws "D27" "21.11"
Set-TextValue $ws "E27" "  -4.69%  "
Set-TextValue $ws "D28" "150.47"
Set-TextValue $ws "E28" "  -0.37%  "
Set-TextValue $ws "D29" "5.247"
Set-TextValue $ws "E29" "  -3.76%  "
Set-TextValue $ws "D30" "132.34"
Set-TextValue $ws "E30" "  -5.09%  "
Set-TextValue $ws "D31" "2.416"
Set-TextValue $ws "E31" "  -2.87%  "
Set-TextValue $ws "D32" "6.946"
Set-TextValue $ws "E32" "  -11.62%  "
Set-TextValue $ws "D33" "1.789.72"
Set-TextValue $ws "E33" "  -2.52%  "
Set-TextValue $ws "D34" "0.9808"
Set-TextValue $ws "E34" "  -6.51%  "
Set-TextValue $ws "D35" "0.07735"
Set-TextValue $ws "E35" "  -4.83%  "
Set-TextValue $ws "D36" "0.02785"
Set-TextValue $ws "E36" "  -6.11%  "
Set-TextValue $ws "D37" "6.286"
Set-TextValue $ws "E37" "  -7.10%  "
Set-TextValue $ws "D38" "0.2552"
Set-TextValue $ws "E38" "  -5.17%  "
Set-TextValue $ws "D39" "10.10"
Set-TextValue $ws "E39" "  -7.75%  "
Set-TextValue $ws "D40" "0.08860"
Set-TextValue $ws "E40" "  -3.33%  "
Set-TextValue $ws "D41" "1.399"
Set-TextValue $ws "E41" "  -2.23%  "
Set-TextValue $ws "D42" "0.7165"
Set-TextValue $ws "E42" "  -5.39%  "
Set-TextValue $ws "D43" "12.76"
Set-TextValue $ws "E43" "  -5.47%  "
Set-TextValue $ws "D44" "15.95"
Set-TextValue $ws "E44" "  -2.53%  "
Set-TextValue $ws "D45" "0.6606"
Set-TextValue $ws "E45" "  -5.09%  "
Set-TextValue $ws "D46" "2.303"
Set-TextValue $ws "E46" "  -6.52%  "
Set-TextValue $ws "D47" "0.9992"
Set-TextValue $ws "E47" "  -0.11%  "
Set-TextValue $ws "E48" "  -2.61%  "
Set-TextValue $ws "D49" "0.08017"
Set-TextValue $ws "E49" "  -3.38%  "
Set-TextValue $ws "D50" "131.19"
Set-TextValue $ws "E50" "  -2.84%  "
Set-TextValue $ws "D51" "1.169"
Set-TextValue $ws "E51" "  -4.08%  "
